$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 80, shifting existing rows 80:100 down to 81:101.
$ws.Rows("80:80").Insert()

# Populate the new row 80 with the new record.
$ws.Range("A80").Value = 4
$ws.Range("B80").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C80").Value = "Los Lagos"
$ws.Range("D80").Value = 45275
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = 300000000
$ws.Range("G80").Value = "Espárragos"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 400
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = 2000
$ws.Range("N80").Value = "$/kilo"
$ws.Range("O80").Value = "Provincia de Linares"
$ws.Range("P80").Value = 2000
$ws.Range("Q80").Value = 1
$ws.Range("R80").Value = "Hortaliza"

# Match the date number-format style used by the rest of column D.
$ws.Range("D80").NumberFormat = $ws.Range("D81").NumberFormat
